$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.177710843373494
$ws.Range("C2").Value = 0.6054216867469879
$ws.Range("J2").Value = 0.003012048192771084
$ws.Range("P2").Value = 0.141566265060241
$ws.Range("S2").Value = 0.07228915662650602
$ws.Range("C3").Value = 0.02415458937198068
$ws.Range("J3").Value = 0.00966183574879227
$ws.Range("P3").Value = 0.8647342995169082
$ws.Range("S3").Value = 0.1014492753623188
$ws.Range("P4").Value = 0.8070175438596491
$ws.Range("S4").Value = 0.1929824561403509
$ws.Range("B6").Value = 0.04824561403508772
$ws.Range("D6").Value = 0.01754385964912281
$ws.Range("F6").Value = 0.04385964912280702
$ws.Range("J6").Value = 0.1885964912280702
$ws.Range("O6").Value = 0.03947368421052631
$ws.Range("Q6").Value = 0.1929824561403509
$ws.Range("R6").Value = 0.1008771929824561
$ws.Range("S6").Value = 0.3684210526315789
$ws.Range("B7").Value = 0.09090909090909091
$ws.Range("D7").Value = 0.0053475935828877
$ws.Range("F7").Value = 0.053475935828877
$ws.Range("J7").Value = 0.1122994652406417
$ws.Range("O7").Value = 0.0106951871657754
$ws.Range("Q7").Value = 0.1925133689839572
$ws.Range("R7").Value = 0.1443850267379679
$ws.Range("S7").Value = 0.3903743315508021
$ws.Range("B8").Value = 0.09677419354838709
$ws.Range("D8").Value = 0.02508960573476703
$ws.Range("F8").Value = 0.06451612903225806
$ws.Range("J8").Value = 0.08243727598566308
$ws.Range("O8").Value = 0.02867383512544803
$ws.Range("Q8").Value = 0.1630824372759857
$ws.Range("R8").Value = 0.1254480286738351
$ws.Range("S8").Value = 0.4139784946236559
$ws.Range("B9").Value = 0.1055276381909548
$ws.Range("D9").Value = 0.02512562814070352
$ws.Range("F9").Value = 0.07537688442211055
$ws.Range("J9").Value = 0.07537688442211055
$ws.Range("O9").Value = 0.02512562814070352
$ws.Range("Q9").Value = 0.185929648241206
$ws.Range("R9").Value = 0.1608040201005025
$ws.Range("S9").Value = 0.3467336683417085
$ws.Range("B10").Value = 0.1166077738515901
$ws.Range("D10").Value = 0.02332155477031802
$ws.Range("E10").Value = 0.001413427561837456
$ws.Range("F10").Value = 0.07208480565371024
$ws.Range("J10").Value = 0.0833922261484099
$ws.Range("O10").Value = 0.0127208480565371
$ws.Range("Q10").Value = 0.2049469964664311
$ws.Range("R10").Value = 0.1293286219081272
$ws.Range("S10").Value = 0.3561837455830389
$ws.Range("G11").Value = 0.1183206106870229
$ws.Range("J11").Value = 0.09541984732824428
$ws.Range("K11").Value = 0.1679389312977099
$ws.Range("L11").Value = 0.6145038167938931
$ws.Range("S11").Value = 0.003816793893129771
$ws.Range("G12").Value = 0.7023809523809523
$ws.Range("J12").Value = 0.2142857142857143
$ws.Range("K12").Value = 0.005952380952380952
$ws.Range("L12").Value = 0.04166666666666666
$ws.Range("S12").Value = 0.03571428571428571
$ws.Range("G13").Value = 0.7916666666666666
$ws.Range("J13").Value = 0.2083333333333333
$ws.Range("F15").Value = 0.003663003663003663
$ws.Range("H15").Value = 0.1684981684981685
$ws.Range("I15").Value = 0.06593406593406594
$ws.Range("J15").Value = 0.3846153846153846
$ws.Range("K15").Value = 0.03296703296703297
$ws.Range("M15").Value = 0.01465201465201465
$ws.Range("N15").Value = 0.007326007326007326
$ws.Range("O15").Value = 0.05128205128205128
$ws.Range("S15").Value = 0.271062271062271
$ws.Range("F16").Value = 0.01520912547528517
$ws.Range("H16").Value = 0.1939163498098859
$ws.Range("I16").Value = 0.06844106463878327
$ws.Range("J16").Value = 0.4372623574144487
$ws.Range("K16").Value = 0.07224334600760456
$ws.Range("M16").Value = 0.01520912547528517
$ws.Range("N16").Value = 0.003802281368821293
$ws.Range("O16").Value = 0.07984790874524715
$ws.Range("S16").Value = 0.1140684410646388
$ws.Range("F17").Value = 0.01622718052738337
$ws.Range("H17").Value = 0.1947261663286004
$ws.Range("I17").Value = 0.08316430020283976
$ws.Range("J17").Value = 0.4482758620689655
$ws.Range("K17").Value = 0.08924949290060852
$ws.Range("M17").Value = 0.01825557809330629
$ws.Range("N17").Value = 0.004056795131845842
$ws.Range("O17").Value = 0.0486815415821501
$ws.Range("S17").Value = 0.0973630831643002
$ws.Range("F18").Value = 0.0060790273556231
$ws.Range("H18").Value = 0.2006079027355623
$ws.Range("I18").Value = 0.0790273556231003
$ws.Range("J18").Value = 0.4437689969604863
$ws.Range("K18").Value = 0.0729483282674772
$ws.Range("M18").Value = 0.02735562310030395
$ws.Range("O18").Value = 0.08206686930091185
$ws.Range("S18").Value = 0.08814589665653495
$ws.Range("F19").Value = 0.00832072617246596
$ws.Range("H19").Value = 0.2239031770045386
$ws.Range("I19").Value = 0.07337367624810892
$ws.Range("J19").Value = 0.4031770045385779
$ws.Range("K19").Value = 0.08850226928895613
$ws.Range("M19").Value = 0.01739788199697428
$ws.Range("N19").Value = 0.0007564296520423601
$ws.Range("O19").Value = 0.08093797276853253
$ws.Range("S19").Value = 0.1036308623298033
